# The sheet originally listed 7 innings (rows 2-8). This edit keeps only the
# "Oct 16 2020 vs Mumbai Indians" innings (previously row 4) as the sole data
# row, moving its values up into row 2, and removes the other innings rows.
# End result: header row (1) + one data row (2) -> dimension A1:K2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with the values that used to live in row 4.
# (D2 "Kolkata Knight Riders" and F2 "Dinesh Karthik <dagger>" are already
# identical across every innings row, so they don't need to be rewritten.)
$ws.Range("A2").Value = " Oct 16 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("E2").Value = "Mumbai Indians"

# G2:K2 hold numeric-looking values that must stay text (matching the
# "numberStoredAsText" cells elsewhere in the sheet), so force a text format
# before assigning, then restore the default style so no stray formatting
# is left behind.
$numRange = $ws.Range("G2:K2")
$numRange.NumberFormat = "@"
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "8"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "50.00"
$numRange.Style = "Normal"

# Drop the other six innings rows (old rows 3-8); their data is no longer
# present anywhere in the sheet, and this shrinks the sheet to A1:K2.
$ws.Rows("3:8").Delete()
